# "Chamada do dia 27" — mark attendance ("C") for the morning and afternoon
# sessions of 27/04 (columns G and H) for every student row (3..49), the
# same way the other date columns were already filled in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3:H49").Value = "C"

# Leave the selection on the range that was just filled in, matching the
# state Excel would show right after the fill.
$ws.Range("G3:H49").Select()
